$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H62").Value = 7918.5
$ws.Range("I62").Value = 0
$ws.Range("J62").Value = 7918.5
$ws.Range("K62").Value = 0
$ws.Range("L62").Value = 7918.5
$ws.Range("M62").ClearContents()
$ws.Range("N62").Value = -9166.5

$ws.Range("H65").Value = 7918.5
$ws.Range("I65").Value = 0
$ws.Range("J65").Value = 7918.5
$ws.Range("K65").Value = 0
$ws.Range("L65").Value = 39592.5
$ws.Range("M65").ClearContents()
$ws.Range("N65").Value = -45832.5

$ws.Range("H105").Value = 39890.332
$ws.Range("J105").Value = 39890.332
$ws.Range("L105").Value = 39890.332
$ws.Range("N105").Value = -46878.332

$ws.Range("H127").Value = 6717
$ws.Range("J127").Value = 4435.5
$ws.Range("L127").Value = 13306.5
$ws.Range("N127").Value = -23226.5

$ws.Range("H132").Value = 942.8333
$ws.Range("I132").Value = 954.93335
$ws.Range("J132").Value = 882.3333
$ws.Range("K132").Value = 2864.80005
$ws.Range("L132").Value = 2646.9999
$ws.Range("M132").Value = -334.8000499999998
$ws.Range("N132").Value = -7706.9999

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H25").Value = 286.25
$ws.Range("I25").Value = 286.25
$ws.Range("K25").Value = 286.25
$ws.Range("M25").Value = 115.75

$ws.Range("H125").Value = 60111
$ws.Range("J125").Value = 60111
$ws.Range("L125").Value = 60111
$ws.Range("N125").Value = -69951

$ws.Range("H132").Value = 2708.4546
$ws.Range("I132").Value = 2399.2856
$ws.Range("K132").Value = 7197.8568
$ws.Range("M132").Value = -4667.8568

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 1980.6666
$ws.Range("I94").Value = 1980.8
$ws.Range("K94").Value = 1980.8
$ws.Range("M94").Value = -1529.8

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H99").Value = 0
$ws.Range("I99").Value = 0
$ws.Range("J99").Value = 0
$ws.Range("K99").Value = 0
$ws.Range("L99").Value = 0
$ws.Range("M99").ClearContents()
$ws.Range("N99").ClearContents()

$ws.Range("H107").Value = 674.17645
$ws.Range("I107").Value = 345.125
$ws.Range("J107").Value = 966.6667
$ws.Range("K107").Value = 345.125
$ws.Range("L107").Value = 966.6667
$ws.Range("M107").Value = 1574.875
$ws.Range("N107").Value = -4806.6667

$ws.Range("H126").Value = 0
$ws.Range("I126").Value = 0
$ws.Range("J126").Value = 0
$ws.Range("K126").Value = 0
$ws.Range("L126").Value = 0
$ws.Range("M126").ClearContents()
$ws.Range("N126").ClearContents()

$ws.Range("H134").Value = 3379.25
$ws.Range("I134").Value = 3379.25
$ws.Range("K134").Value = 10137.75
$ws.Range("M134").Value = -7602.75

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H3").Value = 3208.5
$ws.Range("I3").Value = 3208.5
$ws.Range("K3").Value = 9625.5
$ws.Range("M3").Value = -9513.5

$ws.Range("H69").Value = 2798.2

$ws.Range("H72").Value = 2798.2

$ws.Range("H134").Value = 9052.223
$ws.Range("I134").Value = 1655
$ws.Range("J134").Value = 18298.75
$ws.Range("K134").Value = 4965
$ws.Range("L134").Value = 54896.25
$ws.Range("M134").Value = 105
$ws.Range("N134").Value = -65036.25

$ws.Range("H139").Value = 1390.6
$ws.Range("I139").Value = 1390.6
$ws.Range("K139").Value = 4171.799999999999
$ws.Range("M139").Value = 968.2000000000007

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 7849.3335
$ws.Range("I70").Value = 7889.7
$ws.Range("J70").Value = 7647.5
$ws.Range("K70").Value = 7889.7
$ws.Range("L70").Value = 7647.5
$ws.Range("M70").Value = -7619.7
$ws.Range("N70").Value = -8187.5

$ws.Range("H73").Value = 7849.3335
$ws.Range("I73").Value = 7889.7
$ws.Range("J73").Value = 7647.5
$ws.Range("K73").Value = 7889.7
$ws.Range("L73").Value = 7647.5
$ws.Range("M73").Value = -6953.7
$ws.Range("N73").Value = -9519.5

$ws.Range("H97").Value = 248
$ws.Range("I97").Value = 268.6
$ws.Range("J97").Value = 145
$ws.Range("K97").Value = 268.6
$ws.Range("L97").Value = 145
$ws.Range("M97").Value = 227.4
$ws.Range("N97").Value = -1137

$ws.Range("H122").Value = 3658.0588
$ws.Range("I122").Value = 1737
$ws.Range("K122").Value = 5211
$ws.Range("M122").Value = -2761

$ws.Range("H126").Value = 4126.6665
$ws.Range("I126").Value = 4190
$ws.Range("K126").Value = 12570
$ws.Range("M126").Value = -10100

$ws.Range("H131").Value = 0
$ws.Range("J131").Value = 0
$ws.Range("L131").Value = 0
$ws.Range("N131").ClearContents()

$ws.Range("H132").Value = 0
$ws.Range("I132").Value = 0
$ws.Range("K132").Value = 0
$ws.Range("M132").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 7043.48
$ws.Range("I7").Value = 2017
$ws.Range("K7").Value = 2017
$ws.Range("M7").Value = -1905

$ws.Range("H16").Value = 839.3333
$ws.Range("I16").Value = 549.7778
$ws.Range("J16").Value = 1708
$ws.Range("K16").Value = 549.7778
$ws.Range("L16").Value = 1708
$ws.Range("M16").Value = -379.7778
$ws.Range("N16").Value = -2048

$ws.Range("H40").Value = 4310.875
$ws.Range("J40").Value = 4555.2856
$ws.Range("L40").Value = 4555.2856
$ws.Range("N40").Value = -4827.2856

$ws.Range("H46").Value = 33723.188
$ws.Range("I46").Value = 64446.5
$ws.Range("K46").Value = 64446.5
$ws.Range("M46").Value = -64258.5

$ws.Range("H68").Value = 2549.375
$ws.Range("I68").Value = 1794.5
$ws.Range("J68").Value = 2801
$ws.Range("K68").Value = 1794.5
$ws.Range("L68").Value = 2801
$ws.Range("M68").Value = -1045.5
$ws.Range("N68").Value = -4299

$ws.Range("H71").Value = 2549.375
$ws.Range("I71").Value = 1794.5
$ws.Range("J71").Value = 2801
$ws.Range("K71").Value = 8972.5
$ws.Range("L71").Value = 14005
$ws.Range("M71").Value = -5228.5
$ws.Range("N71").Value = -21493

$ws.Range("H93").Value = 985.44446
$ws.Range("I93").Value = 796
$ws.Range("J93").Value = 1137
$ws.Range("K93").Value = 796
$ws.Range("L93").Value = 1137
$ws.Range("M93").Value = 452
$ws.Range("N93").Value = -3633

$ws.Range("H126").Value = 7043.48
$ws.Range("I126").Value = 2017
$ws.Range("K126").Value = 6051
$ws.Range("M126").Value = -3581

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H41").Value = 0
$ws.Range("J41").Value = 0
$ws.Range("L41").Value = 0
$ws.Range("N41").ClearContents()

$ws.Range("H62").Value = 13846.167
$ws.Range("J62").Value = 13655.5
$ws.Range("L62").Value = 13655.5
$ws.Range("N62").Value = -14903.5

$ws.Range("H65").Value = 13846.167
$ws.Range("J65").Value = 13655.5
$ws.Range("L65").Value = 68277.5
$ws.Range("N65").Value = -74517.5

$ws.Range("H122").Value = 2500
$ws.Range("J122").Value = 2500
$ws.Range("L122").Value = 7500
$ws.Range("N122").Value = -12400
